# Scheduled-runner price refresh: updates the cached market/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) for a batch of leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets. Values below mirror the
# latest pull; a few previously-empty profit cells now get their first value.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1825.5165
$ws.Range("I15").Value = 1825.5165
$ws.Range("K15").Value = 5476.5495
$ws.Range("M15").Value = -5307.5495

$ws.Range("H17").Value = 5669034.5
$ws.Range("J17").Value = 6002395
$ws.Range("L17").Value = 18007185
$ws.Range("N17").Value = -18007521

$ws.Range("H33").Value = 172.4
$ws.Range("I33").Value = 172.4
$ws.Range("K33").Value = 172.4
$ws.Range("M33").Value = 56.59999999999999

$ws.Range("H43").Value = 645.5454999999999
$ws.Range("I43").Value = 396.66666
$ws.Range("J43").Value = 738.875
$ws.Range("K43").Value = 396.66666
$ws.Range("L43").Value = 738.875
$ws.Range("M43").Value = -327.66666
$ws.Range("N43").Value = -876.875

$ws.Range("H93").Value = 28250
$ws.Range("J93").Value = 28250
$ws.Range("L93").Value = 28250
$ws.Range("N93").Value = -33242

$ws.Range("H100").Value = 1878.6086
$ws.Range("I100").Value = 1450.6154
$ws.Range("J100").Value = 2435
$ws.Range("K100").Value = 1450.6154
$ws.Range("L100").Value = 2435
$ws.Range("M100").Value = -909.6153999999999
$ws.Range("N100").Value = -3517

$ws.Range("H103").Value = 41979520
$ws.Range("I103").Value = 83958590
$ws.Range("J103").Value = 450
$ws.Range("K103").Value = 251875770
$ws.Range("L103").Value = 1350
$ws.Range("M103").Value = -251875184
$ws.Range("N103").Value = -2522

$ws.Range("H116").Value = 4154.091
$ws.Range("J116").Value = 4686.5713
$ws.Range("L116").Value = 4686.5713
$ws.Range("N116").Value = -11570.5713

$ws.Range("H129").Value = 176429.14
$ws.Range("J129").Value = 186195.58
$ws.Range("L129").Value = 558586.74
$ws.Range("N129").Value = -568586.74

$ws.Range("H140").Value = 50722.855
$ws.Range("J140").Value = 50722.855
$ws.Range("L140").Value = 50722.855
$ws.Range("N140").Value = -61082.855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5769.5366
$ws.Range("I32").Value = 4527.3584
$ws.Range("J32").Value = 11317.934
$ws.Range("K32").Value = 4527.3584
$ws.Range("L32").Value = 11317.934
$ws.Range("M32").Value = -4240.3584
$ws.Range("N32").Value = -11891.934

$ws.Range("H45").Value = 2565.4194
$ws.Range("I45").Value = 2136.3572
$ws.Range("J45").Value = 2918.7646
$ws.Range("K45").Value = 2136.3572
$ws.Range("L45").Value = 2918.7646
$ws.Range("M45").Value = -1759.3572
$ws.Range("N45").Value = -3672.7646

$ws.Range("H88").Value = 144535.58
$ws.Range("J88").Value = 334484.34
$ws.Range("L88").Value = 334484.34
$ws.Range("N88").Value = -335296.34

$ws.Range("H91").Value = 144535.58
$ws.Range("J91").Value = 334484.34
$ws.Range("L91").Value = 334484.34
$ws.Range("N91").Value = -337292.34

$ws.Range("H110").Value = 954.3333
$ws.Range("I110").Value = 898.625
$ws.Range("K110").Value = 898.625
$ws.Range("M110").Value = 1146.375

$ws.Range("H132").Value = 28806.3
$ws.Range("I132").Value = 3954.923
$ws.Range("K132").Value = 11864.769
$ws.Range("M132").Value = -9334.769

$ws.Range("H137").Value = 42474.6
$ws.Range("J137").Value = 42343.25
$ws.Range("L137").Value = 42343.25
$ws.Range("N137").Value = -52543.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2412.261
$ws.Range("I20").Value = 2745.6667
$ws.Range("K20").Value = 2745.6667
$ws.Range("M20").Value = -2498.6667

$ws.Range("H94").Value = 1032.0769
$ws.Range("I94").Value = 837.125
$ws.Range("K94").Value = 837.125
$ws.Range("M94").Value = -386.125

$ws.Range("H107").Value = 1041.1
$ws.Range("I107").Value = 926.375
$ws.Range("K107").Value = 926.375
$ws.Range("M107").Value = 993.625

$ws.Range("H134").Value = 5240.2383
$ws.Range("I134").Value = 6203.357
$ws.Range("K134").Value = 18610.071
$ws.Range("M134").Value = -16075.071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3370.6
$ws.Range("I31").Value = 1730.1538
$ws.Range("K31").Value = 1730.1538
$ws.Range("M31").Value = -1435.1538

$ws.Range("H34").Value = 3370.6
$ws.Range("I34").Value = 1730.1538
$ws.Range("K34").Value = 1730.1538
$ws.Range("M34").Value = -1528.1538

$ws.Range("H86").Value = 16506.691
$ws.Range("I86").Value = 8314.25
$ws.Range("K86").Value = 8314.25
$ws.Range("M86").Value = -7191.25

$ws.Range("H89").Value = 16506.691
$ws.Range("I89").Value = 8314.25
$ws.Range("K89").Value = 41571.25
$ws.Range("M89").Value = -35955.25

$ws.Range("H132").Value = 4948.923
$ws.Range("I132").Value = 3669.3333
$ws.Range("J132").Value = 7828
$ws.Range("K132").Value = 11007.9999
$ws.Range("L132").Value = 23484
$ws.Range("M132").Value = -8477.999899999999
$ws.Range("N132").Value = -28544

$ws.Range("H134").Value = 1716.6666
$ws.Range("I134").Value = 1433.3334
$ws.Range("K134").Value = 4300.0002
$ws.Range("M134").Value = -1765.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4730.3184
$ws.Range("I107").Value = 5171
$ws.Range("J107").Value = 323.5
$ws.Range("K107").Value = 15513
$ws.Range("L107").Value = 970.5
$ws.Range("M107").Value = -13593
$ws.Range("N107").Value = -4810.5

$ws.Range("H113").Value = 757.6842
$ws.Range("I113").Value = 626.1
$ws.Range("J113").Value = 903.8889
$ws.Range("K113").Value = 1878.3
$ws.Range("L113").Value = 2711.6667
$ws.Range("M113").Value = 291.6999999999998
$ws.Range("N113").Value = -7051.6667

$ws.Range("H131").Value = 702.12
$ws.Range("J131").Value = 707.1919
$ws.Range("L131").Value = 2121.5757
$ws.Range("N131").Value = -12201.5757

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 899
$ws.Range("J29").Value = 899
$ws.Range("L29").Value = 899
$ws.Range("N29").Value = -1479

$ws.Range("H80").Value = 3436.182
$ws.Range("I80").Value = 2837.5
$ws.Range("J80").Value = 3778.2856
$ws.Range("K80").Value = 2837.5
$ws.Range("L80").Value = 3778.2856
$ws.Range("M80").Value = -1839.5
$ws.Range("N80").Value = -5774.2856

$ws.Range("H83").Value = 3436.182
$ws.Range("I83").Value = 2837.5
$ws.Range("J83").Value = 3778.2856
$ws.Range("K83").Value = 14187.5
$ws.Range("L83").Value = 18891.428
$ws.Range("M83").Value = -9195.5
$ws.Range("N83").Value = -28875.428

$ws.Range("H102").Value = 1480.6061
$ws.Range("I102").Value = 1225.4642
$ws.Range("J102").Value = 2909.4
$ws.Range("K102").Value = 1225.4642
$ws.Range("L102").Value = 2909.4
$ws.Range("M102").Value = 396.5358000000001
$ws.Range("N102").Value = -6153.4

$ws.Range("H132").Value = 34984.41
$ws.Range("I132").Value = 6210.5386
$ws.Range("K132").Value = 18631.6158
$ws.Range("M132").Value = -16101.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3218.8125
$ws.Range("I22").Value = 4222.8184
$ws.Range("J22").Value = 1010
$ws.Range("K22").Value = 4222.8184
$ws.Range("L22").Value = 1010
$ws.Range("M22").Value = -3927.8184
$ws.Range("N22").Value = -1600

$ws.Range("H27").Value = 3218.8125
$ws.Range("I27").Value = 4222.8184
$ws.Range("J27").Value = 1010
$ws.Range("K27").Value = 4222.8184
$ws.Range("L27").Value = 1010
$ws.Range("M27").Value = -4115.8184
$ws.Range("N27").Value = -1224

$ws.Range("H29").Value = 50003000
$ws.Range("I29").Value = 100000000
$ws.Range("J29").Value = 6000
$ws.Range("K29").Value = 100000000
$ws.Range("L29").Value = 6000
$ws.Range("M29").Value = -99999705
$ws.Range("N29").Value = -6590

$ws.Range("H40").Value = 3387.5
$ws.Range("I40").Value = 2926.7856
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2926.7856
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2790.7856
$ws.Range("N40").Value = -5272

$ws.Range("H46").Value = 2554.4546
$ws.Range("I46").Value = 3049.875
$ws.Range("K46").Value = 3049.875
$ws.Range("M46").Value = -2861.875

$ws.Range("H55").Value = 913.63635
$ws.Range("I55").Value = 913.63635
$ws.Range("K55").Value = 913.63635
$ws.Range("M55").Value = -740.63635

$ws.Range("H82").Value = 939.875
$ws.Range("I82").Value = 902.7143
$ws.Range("K82").Value = 902.7143
$ws.Range("M82").Value = -541.7143

$ws.Range("H85").Value = 939.875
$ws.Range("I85").Value = 902.7143
$ws.Range("K85").Value = 902.7143
$ws.Range("M85").Value = 345.2857

$ws.Range("H122").Value = 1511796.2
$ws.Range("I122").Value = 1963605.2
$ws.Range("K122").Value = 5890815.6
$ws.Range("M122").Value = -5888365.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1528.25
$ws.Range("I122").Value = 1482.8422
$ws.Range("J122").Value = 1700.8
$ws.Range("K122").Value = 4448.5266
$ws.Range("L122").Value = 5102.4
$ws.Range("M122").Value = -1998.5266
$ws.Range("N122").Value = -10002.4

$ws.Range("H141").Value = 56333.332
$ws.Range("J141").Value = 56333.332
$ws.Range("L141").Value = 56333.332
$ws.Range("N141").Value = -66693.33199999999
